# Update gh-pages output data (column F = "想去人数") for sheets "展览" and "全部类型"

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1143
$ws1.Range("F3").Value = 413
$ws1.Range("F4").Value = 257
$ws1.Range("F6").Value = 6
$ws1.Range("F7").Value = 12145
$ws1.Range("F9").Value = 8
$ws1.Range("F10").Value = 108
$ws1.Range("F11").Value = 11916
$ws1.Range("F12").Value = 4780
$ws1.Range("F13").Value = 605
$ws1.Range("F14").Value = 93
$ws1.Range("F15").Value = 34
$ws1.Range("F18").Value = 937
$ws1.Range("F21").Value = 60

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1143
$ws4.Range("F3").Value = 413
$ws4.Range("F4").Value = 257
$ws4.Range("F8").Value = 6
$ws4.Range("F9").Value = 12145
$ws4.Range("F11").Value = 8
$ws4.Range("F12").Value = 108
$ws4.Range("F13").Value = 11916
$ws4.Range("F14").Value = 4780
$ws4.Range("F15").Value = 605
$ws4.Range("F16").Value = 93
$ws4.Range("F17").Value = 34
$ws4.Range("F20").Value = 937
$ws4.Range("F23").Value = 60
